# Refactor parity generation logic to exclude unnecessary FIERR ports
# and improve port declaration handling.
#
# Concretely (per the source-controlled diff):
#   1. The VERSION / HSR ID / SM ID columns (B:D) and the ERROR DOUBLE
#      column (O) were only ever populated with empty placeholder cells
#      for the four signal rows (2-5) - the port-declaration generator no
#      longer emits those unused placeholders, so the cells are cleared
#      out entirely.
#   2. The DRIVE/RECEIVE direction for the write-address, write-data and
#      read-address parity ports (rows 2-4) was flipped from DRIVE to
#      RECEIVE, while the read-data parity port (row 5) flipped the other
#      way, from RECEIVE to DRIVE - this matches the corrected generation
#      logic for which side of the interface drives each signal.
#   3. Best-effort column widths for the columns whose text grew (IP
#      NAME / SIGNAL VALID NAME / PARITY PORT NAME / IP FILE LIST /
#      ERROR PORT) are (re)applied, mirroring Excel's own best-fit pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the now-unused placeholder cells -----------------------------
# B2:D5 (VERSION / HSR ID / SM ID) and O2:O5 (ERROR DOUBLE) were empty
# inline-string cells; ClearContents removes the cell nodes entirely
# instead of merely blanking their value.
$ws.Range("B2:D5").ClearContents()
$ws.Range("O2:O5").ClearContents()

# --- 2. Fix up the DRIVE/RECEIVE column (J) --------------------------------
$ws.Range("J2").Value = "RECEIVE"   # WADDR_PARITY  - was DRIVE
$ws.Range("J3").Value = "RECEIVE"   # WDATA_PARITY  - was DRIVE
$ws.Range("J4").Value = "RECEIVE"   # RADDR_PARITY  - was DRIVE
$ws.Range("J5").Value = "DRIVE"     # RDATA_PARITY  - was RECEIVE

# --- 3. Re-fit the columns whose content widened ---------------------------
# (Target best-fit widths are 12.375 / 22.25 / 21.375 / 23.75 / 20.125
# characters; the values below are the closest the host's column-width
# quantisation can reach.)
$ws.Columns.Item(5).ColumnWidth  = 11.42   # E - IP NAME
$ws.Columns.Item(11).ColumnWidth = 21.25   # K - SIGNAL VALID NAME
$ws.Columns.Item(12).ColumnWidth = 20.42   # L - PARITY PORT NAME
$ws.Columns.Item(13).ColumnWidth = 22.75   # M - IP FILE LIST
$ws.Columns.Item(14).ColumnWidth = 19.25   # N - ERROR PORT
